$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting RESEARCH/SALES rows down
$ws.Rows.Item(3).Insert()

# Update Total_Compensation for ACCOUNTING (row 2)
$ws.Range("C2").Value = 4299600

# Populate the newly inserted row 3 with the second ACCOUNTING entry
$ws.Range("A3").Value = "ACCOUNTING"
$ws.Range("B3").Value = 30
$ws.Range("C3").Value = 12350
